# Applies the crypto-dashboard refresh described in the commit:
#   "Updated cryptos list on Wed Jan  3 04:50:35 UTC 2024 with GitHub Actions"
#
# Plain .Value assignment lets the COM layer auto-detect numbers, which
# would silently coerce price strings like "40.80" or "8.48" into the
# numeric 40.8 / 8.48 (losing the trailing zero / becoming t="n" instead
# of the original inline-string cell). Price cells that parse as plain
# floats are therefore written via Set-TextValue, which flips the cell to
# the text format ("@") just long enough to force string storage, then
# calls ClearFormats() to drop that temporary number format again so the
# cell keeps its original (default) style index.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$ws.Range("D2").Value = "45.346.47"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "2.369.16"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  -0.09%  "
Set-TextValue $ws.Range("D5") "314.11"
$ws.Range("E5").Value = "  -1.08%  "
Set-TextValue $ws.Range("D6") "107.95"
$ws.Range("E6").Value = "  -3.54%  "
Set-TextValue $ws.Range("D7") "0.633"
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("E9").Value = "  -2.81%  "
Set-TextValue $ws.Range("D10") "40.80"
$ws.Range("E10").Value = "  -3.67%  "
Set-TextValue $ws.Range("D11") "0.0917"
$ws.Range("E11").Value = "  -1.52%  "
Set-TextValue $ws.Range("D12") "8.48"
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("E14").Value = "  -3.90%  "
$ws.Range("D15").Value = "2.728.88"
$ws.Range("E15").Value = "  -0.64%  "
Set-TextValue $ws.Range("D16") "15.33"
$ws.Range("E16").Value = "  -3.39%  "
$ws.Range("D17").Value = "2.365.60"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").Value = "45.292.30"
$ws.Range("E18").Value = "  -0.11%  "
Set-TextValue $ws.Range("D19") "15.44"
$ws.Range("E19").Value = "  +17.18%  "
$ws.Range("E20").Value = "  -5.33%  "
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("E22").Value = "  +1.84%  "
Set-TextValue $ws.Range("D23") "73.30"
$ws.Range("E23").Value = "  -2.38%  "
Set-TextValue $ws.Range("D24") "261.04"
$ws.Range("E24").Value = "  -3.19%  "
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D27") "11.15"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D28") "7.45"
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("E29").Value = "  -1.68%  "
Set-TextValue $ws.Range("D30") "0.0967"
$ws.Range("E30").Value = "  +2.49%  "
Set-TextValue $ws.Range("D31") "22.31"
$ws.Range("E31").Value = "  -2.56%  "
Set-TextValue $ws.Range("D32") "37.06"
$ws.Range("E32").Value = "  -3.84%  "
Set-TextValue $ws.Range("D33") "166.43"
Set-TextValue $ws.Range("D34") "2.87"
$ws.Range("E34").Value = "  -4.54%  "
$ws.Range("E35").Value = "  -2.03%  "
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("E37").Value = "  -4.00%  "
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("E39").Value = "  +7.62%  "
$ws.Range("E40").Value = "  -6.28%  "
$ws.Range("E41").Value = "  -3.06%  "
Set-TextValue $ws.Range("D42") "98.52"
$ws.Range("E42").Value = "  -6.32%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue $ws.Range("D43") "70.13"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D44") "13.12"
$ws.Range("E44").Value = "  -3.18%  "
$ws.Range("E45").Value = "  -5.24%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("E47").Value = "  +2.68%  "
$ws.Range("D48").Value = "1.819.13"
$ws.Range("E48").Value = "  +10.57%  "
Set-TextValue $ws.Range("D49") "82.84"
$ws.Range("E49").Value = "  +4.77%  "
Set-TextValue $ws.Range("D50") "111.16"
$ws.Range("E50").Value = "  -5.71%  "
$ws.Range("E51").Value = "  -2.07%  "
